# Update the "想去人数" (want-to-go count) figures for two events that
# appear on both the "展览" sheet and the aggregated "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 359
    $ws.Range("F5").Value = 296
}
